$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 1.399810963551616
$ws.Range("F1").Value = -1.570796390562869

$ws.Range("E2").Value = 1.398771000953773
$ws.Range("F2").Value = -1.57079639084348

$ws.Range("E3").Value = 1.394111279037471
$ws.Range("F3").Value = -1.570796392100804

$ws.Range("E4").Value = 1.387555492716261
$ws.Range("F4").Value = -1.57079639386974

$ws.Range("E5").Value = 1.382895770799958
$ws.Range("F5").Value = -1.570796395127064

$ws.Range("E6").Value = 1.381855808202115
$ws.Range("F6").Value = -1.570796395407675
